$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.045.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.261.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.57"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.60%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.450"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +14.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +19.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.601.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.836"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.265.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.880.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.74%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.00%  "
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +23.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.137"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0681"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.73"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.10%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.40"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.99%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0963"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("B46").Value = "TerraClassic"
$ws.Range("C46").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000213"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.34%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +17.76%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.446.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.82%  "
